# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handback DateTime"
# timestamps for the 597be984-... / bbda21e6-... handback entries.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 597be984-...md
$overview.Range("G4").Value = "2016-08-19 00:43:58"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 597be984-...
$zhcn.Range("H4").Value = "2016-08-19 00:43:53"
$zhcn.Range("K4").Value = "2016-08-19 00:44:15"

# de-de sheet: Correspond Handback DateTime for 597be984-...
$dede.Range("K4").Value = "2016-08-19 00:44:23"
